$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Trophy" -> "Trophic level" in the Internal toolbox name column (D10)
$ws.Range("D10").Value = "Trophic level"

# Update the accompanying comment text in E10 to match the new name
$ws.Range("E10").Value = "Trophic level will automatically be put on taxon level if size classes are equal. Are different for Unicell etc."

# Move the active selection to D46, as captured in the saved view state
$ws.Range("D46").Select()
